$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Drop the "Email" column's mailto: hyperlinks (both the link objects and
# the left-over hyperlink character style they leave behind on the cells).
# ---------------------------------------------------------------------
$ws.Hyperlinks.Delete()
$ws.Range("C2").Style = "Normal"
$ws.Range("C3").Style = "Normal"
$ws.Range("C2:C3").HorizontalAlignment = -4108   # xlCenter
$wb.Styles("Hyperlink").Delete()

# ---------------------------------------------------------------------
# Re-derive the two new header-row look-alikes from the existing bordered
# header style before we start moving values around:
#   - E1 keeps the bold header font but loses every border (it is now a
#     trailing blank header cell).
#   - C1 (the new "Password" header) keeps a border on three sides only
#     (its right edge is left open).
# ---------------------------------------------------------------------
$ws.Range("A1").Copy($ws.Range("E1"))
$ws.Range("E1").Borders.LineStyle = -4142           # drop all four edges in one go
$ws.Range("E1").ClearContents()

$ws.Range("A1").Copy($ws.Range("C1"))
$ws.Range("C1").Borders.Item(10).LineStyle = -4142  # xlEdgeRight -> none

# ---------------------------------------------------------------------
# Column C: was "Email", now becomes "Password".
# ---------------------------------------------------------------------
$ws.Range("C1").Value = "Password"
$ws.Range("C2").Value = "sample123"
$ws.Range("C3").Value = "sample123"

# ---------------------------------------------------------------------
# Column D: was "Password"/sample passwords, now becomes "Aktif"
# (the active-flag column that used to live in E).
# ---------------------------------------------------------------------
$ws.Range("D1").Value = "Aktif"
$ws.Range("D2").Value = 1
$ws.Range("D3").Value = 0

# ---------------------------------------------------------------------
# Column E: used to hold the active-flag values, now just a blank,
# still-styled trailing column.
# ---------------------------------------------------------------------
$ws.Range("E2:E3").ClearContents()

# ---------------------------------------------------------------------
# Column F ("Roles") is removed entirely - content and formatting.
# ---------------------------------------------------------------------
$ws.Range("F1:F3").Clear()

# ---------------------------------------------------------------------
# View: zoomed in a bit, selection parked on the sample password cell.
# ---------------------------------------------------------------------
$ws.Application.ActiveWindow.Zoom = 130
$ws.Range("C2").Select()
